$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1. "...We modified the script and used some options of imagemagick..."
#    -> "...We modified the script and explored some options of imagemagick..."
# -----------------------------------------------------------------
$d.Content.Find.Execute("We modified the script and used some options of", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "We modified the script and explored some options of", 2) | Out-Null

# -----------------------------------------------------------------
# 2. Bullet about keywords: append '  OR "AERIAL"' right after the
#    closing quote of ...UNIDENTIFIED FLYING OBJECT"
# -----------------------------------------------------------------
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("REPORT OF UNIDENTIFIED FLYING OBJECT”", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $rng2.Collapse(0)
    $rng2.InsertAfter("  OR “AERIAL”")
}

# -----------------------------------------------------------------
# 3. "...contained the keywords "FLYING OBJECT". We ignored..."
#    -> "...contained the keywords "FLYING OBJECT" or "AERIAL". We ignored..."
# -----------------------------------------------------------------
$rng3 = $d.Content
$found3 = $rng3.Find.Execute("pdf files contained the keywords “FLYING OBJECT”", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found3) {
    $rng3.Collapse(0)
    $rng3.InsertAfter(" or “AERIAL”")
}
